# Auto-generated edit script: reorders rows 7-11 and 45-49 content
# to match the target diff (cells are swapped/rotated between rows;
# row numbers themselves do not change, only their cell contents).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 ---
$ws.Range("A7").Value2 = 130752192
$ws.Range("B7").Value2 = 57881
$ws.Range("D7").Value2 = 'NT'
$ws.Range("E7").Value2 = 100049
$ws.Range("F7").Value2 = 'Spillkråka'
$ws.Range("G7").Value2 = 'Dryocopus martius'
$ws.Range("H7").Value2 = '(Linnaeus, 1758)'
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").Value2 = 'äldre spår'
$ws.Range("N7").ClearContents()
$ws.Range("P7").Value2 = 'Truppan, Dlr'
$ws.Range("Q7").Value2 = 490715
$ws.Range("R7").Value2 = 6763290
$ws.Range("S7").Value2 = 10
$ws.Range("T7").Value2 = 'Dalarna'
$ws.Range("U7").Value2 = 'Mora'
$ws.Range("V7").Value2 = 'Dalarna'
$ws.Range("W7").Value2 = 'Mora'
$ws.Range("Z7").Value2 = '11:43'
$ws.Range("AB7").Value2 = '11:43'
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").Value2 = $False
$ws.Range("AE7").Value2 = $False
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").Value2 = $False
$ws.Range("AT7").ClearContents()
$ws.Range("AW7").Value2 = 'Håkan Thenander'
$ws.Range("AX7").Value2 = 'Håkan Thenander, Bo karlstens'
$ws.Range("AY7").ClearContents()

# --- Row 8 ---
$ws.Range("A8").Value2 = 130789501
$ws.Range("B8").Value2 = 79243
$ws.Range("D8").Value2 = 'NT'
$ws.Range("E8").Value2 = 6425
$ws.Range("F8").Value2 = 'Garnlav'
$ws.Range("G8").Value2 = 'Alectoria sarmentosa'
$ws.Range("H8").Value2 = '(Ach.) Ach.'
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").Value2 = 'Kråkbackarna, Dlr'
$ws.Range("Q8").Value2 = 490713
$ws.Range("R8").Value2 = 6763507
$ws.Range("S8").Value2 = 10
$ws.Range("T8").Value2 = 'Dalarna'
$ws.Range("U8").Value2 = 'Mora'
$ws.Range("V8").Value2 = 'Dalarna'
$ws.Range("W8").Value2 = 'Mora'
$ws.Range("Z8").Value2 = '12:36'
$ws.Range("AB8").Value2 = '12:36'
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").Value2 = $False
$ws.Range("AE8").Value2 = $False
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").Value2 = $False
$ws.Range("AT8").ClearContents()
$ws.Range("AW8").Value2 = 'Bo karlstens'
$ws.Range("AX8").Value2 = 'Bo karlstens, Håkan Thenander'
$ws.Range("AY8").ClearContents()

# --- Row 9 ---
$ws.Range("A9").Value2 = 130789509
$ws.Range("B9").Value2 = 79243
$ws.Range("D9").Value2 = 'NT'
$ws.Range("E9").Value2 = 6425
$ws.Range("F9").Value2 = 'Garnlav'
$ws.Range("G9").Value2 = 'Alectoria sarmentosa'
$ws.Range("H9").Value2 = '(Ach.) Ach.'
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").Value2 = 'Kråkbackarna, Dlr'
$ws.Range("Q9").Value2 = 490693
$ws.Range("R9").Value2 = 6763417
$ws.Range("S9").Value2 = 10
$ws.Range("T9").Value2 = 'Dalarna'
$ws.Range("U9").Value2 = 'Mora'
$ws.Range("V9").Value2 = 'Dalarna'
$ws.Range("W9").Value2 = 'Mora'
$ws.Range("Z9").Value2 = '12:05'
$ws.Range("AB9").Value2 = '12:05'
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").Value2 = $False
$ws.Range("AE9").Value2 = $False
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").Value2 = $False
$ws.Range("AT9").ClearContents()
$ws.Range("AW9").Value2 = 'Bo karlstens'
$ws.Range("AX9").Value2 = 'Bo karlstens, Håkan Thenander'
$ws.Range("AY9").ClearContents()

# --- Row 10 ---
$ws.Range("A10").Value2 = 130752527
$ws.Range("B10").Value2 = 79243
$ws.Range("D10").Value2 = 'NT'
$ws.Range("E10").Value2 = 6425
$ws.Range("F10").Value2 = 'Garnlav'
$ws.Range("G10").Value2 = 'Alectoria sarmentosa'
$ws.Range("H10").Value2 = '(Ach.) Ach.'
$ws.Range("I10").ClearContents()
$ws.Range("K10").ClearContents()
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("P10").Value2 = 'Kråkbackarna, Dlr'
$ws.Range("Q10").Value2 = 490673
$ws.Range("R10").Value2 = 6763435
$ws.Range("S10").Value2 = 10
$ws.Range("T10").Value2 = 'Dalarna'
$ws.Range("U10").Value2 = 'Mora'
$ws.Range("V10").Value2 = 'Dalarna'
$ws.Range("W10").Value2 = 'Mora'
$ws.Range("Z10").Value2 = '11:43'
$ws.Range("AB10").Value2 = '11:43'
$ws.Range("AC10").ClearContents()
$ws.Range("AD10").Value2 = $False
$ws.Range("AE10").Value2 = $False
$ws.Range("AF10").ClearContents()
$ws.Range("AG10").Value2 = $False
$ws.Range("AT10").ClearContents()
$ws.Range("AW10").Value2 = 'Håkan Thenander'
$ws.Range("AX10").Value2 = 'Håkan Thenander, Bo karlstens'
$ws.Range("AY10").ClearContents()

# --- Row 11 ---
$ws.Range("A11").Value2 = 130789462
$ws.Range("B11").Value2 = 57884
$ws.Range("D11").Value2 = 'NT'
$ws.Range("E11").Value2 = 100109
$ws.Range("F11").Value2 = 'Tretåig hackspett'
$ws.Range("G11").Value2 = 'Picoides tridactylus'
$ws.Range("H11").Value2 = '(Linnaeus, 1758)'
$ws.Range("I11").ClearContents()
$ws.Range("K11").ClearContents()
$ws.Range("L11").ClearContents()
$ws.Range("M11").Value2 = 'äldre spår'
$ws.Range("N11").ClearContents()
$ws.Range("P11").Value2 = 'Kråkbackarna, Dlr'
$ws.Range("Q11").Value2 = 490705
$ws.Range("R11").Value2 = 6763439
$ws.Range("S11").Value2 = 10
$ws.Range("T11").Value2 = 'Dalarna'
$ws.Range("U11").Value2 = 'Mora'
$ws.Range("V11").Value2 = 'Dalarna'
$ws.Range("W11").Value2 = 'Mora'
$ws.Range("Z11").Value2 = '12:09'
$ws.Range("AB11").Value2 = '12:09'
$ws.Range("AC11").Value2 = 'Äldre ring hack på tall'
$ws.Range("AD11").Value2 = $False
$ws.Range("AE11").Value2 = $False
$ws.Range("AF11").ClearContents()
$ws.Range("AG11").Value2 = $False
$ws.Range("AT11").ClearContents()
$ws.Range("AW11").Value2 = 'Bo karlstens'
$ws.Range("AX11").Value2 = 'Bo karlstens, Håkan Thenander'
$ws.Range("AY11").ClearContents()

# --- Row 45 ---
$ws.Range("A45").Value2 = 130807362
$ws.Range("B45").Value2 = 79243
$ws.Range("D45").Value2 = 'NT'
$ws.Range("E45").Value2 = 6425
$ws.Range("F45").Value2 = 'Garnlav'
$ws.Range("G45").Value2 = 'Alectoria sarmentosa'
$ws.Range("H45").Value2 = '(Ach.) Ach.'
$ws.Range("I45").ClearContents()
$ws.Range("K45").ClearContents()
$ws.Range("L45").ClearContents()
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()
$ws.Range("P45").Value2 = 'Truppan, Dlr'
$ws.Range("Q45").Value2 = 490990
$ws.Range("R45").Value2 = 6763195
$ws.Range("S45").Value2 = 10
$ws.Range("T45").Value2 = 'Dalarna'
$ws.Range("U45").Value2 = 'Mora'
$ws.Range("V45").Value2 = 'Dalarna'
$ws.Range("W45").Value2 = 'Mora'
$ws.Range("Z45").Value2 = '11:46'
$ws.Range("AB45").Value2 = '11:46'
$ws.Range("AC45").ClearContents()
$ws.Range("AD45").Value2 = $False
$ws.Range("AE45").Value2 = $False
$ws.Range("AF45").ClearContents()
$ws.Range("AG45").Value2 = $False
$ws.Range("AT45").ClearContents()
$ws.Range("AW45").Value2 = 'Håkan Thenander'
$ws.Range("AX45").Value2 = 'Håkan Thenander'
$ws.Range("AY45").ClearContents()

# --- Row 46 ---
$ws.Range("A46").Value2 = 130815959
$ws.Range("B46").Value2 = 57884
$ws.Range("D46").Value2 = 'NT'
$ws.Range("E46").Value2 = 100109
$ws.Range("F46").Value2 = 'Tretåig hackspett'
$ws.Range("G46").Value2 = 'Picoides tridactylus'
$ws.Range("H46").Value2 = '(Linnaeus, 1758)'
$ws.Range("I46").ClearContents()
$ws.Range("K46").ClearContents()
$ws.Range("L46").ClearContents()
$ws.Range("M46").Value2 = 'äldre spår'
$ws.Range("N46").ClearContents()
$ws.Range("P46").Value2 = 'Kråkbackarna, Dlr'
$ws.Range("Q46").Value2 = 490727
$ws.Range("R46").Value2 = 6763422
$ws.Range("S46").Value2 = 10
$ws.Range("T46").Value2 = 'Dalarna'
$ws.Range("U46").Value2 = 'Mora'
$ws.Range("V46").Value2 = 'Dalarna'
$ws.Range("W46").Value2 = 'Mora'
$ws.Range("Z46").Value2 = '11:46'
$ws.Range("AB46").Value2 = '11:46'
$ws.Range("AC46").Value2 = '3 bilder på tall'
$ws.Range("AD46").Value2 = $False
$ws.Range("AE46").Value2 = $False
$ws.Range("AF46").ClearContents()
$ws.Range("AG46").Value2 = $False
$ws.Range("AT46").ClearContents()
$ws.Range("AW46").Value2 = 'Håkan Thenander'
$ws.Range("AX46").Value2 = 'Håkan Thenander'
$ws.Range("AY46").ClearContents()

# --- Row 48 ---
$ws.Range("A48").Value2 = 130814260
$ws.Range("B48").Value2 = 57884
$ws.Range("D48").Value2 = 'NT'
$ws.Range("E48").Value2 = 100109
$ws.Range("F48").Value2 = 'Tretåig hackspett'
$ws.Range("G48").Value2 = 'Picoides tridactylus'
$ws.Range("H48").Value2 = '(Linnaeus, 1758)'
$ws.Range("I48").ClearContents()
$ws.Range("K48").ClearContents()
$ws.Range("L48").ClearContents()
$ws.Range("M48").Value2 = 'färska spår'
$ws.Range("N48").ClearContents()
$ws.Range("P48").Value2 = 'Truppan, Dlr'
$ws.Range("Q48").Value2 = 490982
$ws.Range("R48").Value2 = 6763304
$ws.Range("S48").Value2 = 10
$ws.Range("T48").Value2 = 'Dalarna'
$ws.Range("U48").Value2 = 'Mora'
$ws.Range("V48").Value2 = 'Dalarna'
$ws.Range("W48").Value2 = 'Mora'
$ws.Range("Z48").Value2 = '11:46'
$ws.Range("AB48").Value2 = '11:46'
$ws.Range("AC48").Value2 = '4 bilder på tall'
$ws.Range("AD48").Value2 = $False
$ws.Range("AE48").Value2 = $False
$ws.Range("AF48").ClearContents()
$ws.Range("AG48").Value2 = $False
$ws.Range("AT48").ClearContents()
$ws.Range("AW48").Value2 = 'Håkan Thenander'
$ws.Range("AX48").Value2 = 'Håkan Thenander'
$ws.Range("AY48").ClearContents()

# --- Row 49 ---
$ws.Range("A49").Value2 = 130815999
$ws.Range("B49").Value2 = 79243
$ws.Range("D49").Value2 = 'NT'
$ws.Range("E49").Value2 = 6425
$ws.Range("F49").Value2 = 'Garnlav'
$ws.Range("G49").Value2 = 'Alectoria sarmentosa'
$ws.Range("H49").Value2 = '(Ach.) Ach.'
$ws.Range("I49").ClearContents()
$ws.Range("K49").ClearContents()
$ws.Range("L49").ClearContents()
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()
$ws.Range("P49").Value2 = 'Kråkbackarna, Dlr'
$ws.Range("Q49").Value2 = 490722
$ws.Range("R49").Value2 = 6763404
$ws.Range("S49").Value2 = 50
$ws.Range("T49").Value2 = 'Dalarna'
$ws.Range("U49").Value2 = 'Mora'
$ws.Range("V49").Value2 = 'Dalarna'
$ws.Range("W49").Value2 = 'Mora'
$ws.Range("Z49").Value2 = '11:46'
$ws.Range("AB49").Value2 = '11:46'
$ws.Range("AC49").Value2 = 'Måttligt i en radie av ca 50 meter. 2 bilder tall'
$ws.Range("AD49").Value2 = $False
$ws.Range("AE49").Value2 = $False
$ws.Range("AF49").ClearContents()
$ws.Range("AG49").Value2 = $False
$ws.Range("AT49").ClearContents()
$ws.Range("AW49").Value2 = 'Håkan Thenander'
$ws.Range("AX49").Value2 = 'Håkan Thenander'
$ws.Range("AY49").ClearContents()
